$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "35.625.89"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.98%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.895.66"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.45%  "
$ws.Range("E4").Value = "  -0.71%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "247.55"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -3.19%  "
$ws.Range("E6").Value = "  -5.19%  "
$ws.Range("E7").Value = "  -0.77%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "43.92"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +8.23%  "
$ws.Range("E9").Value = "  -4.47%  "
$ws.Range("E10").Value = "  -2.20%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0970"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.87%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "13.12"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.59%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.170.67"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.43%  "
$ws.Range("E14").Value = "  +0.73%  "
$ws.Range("B15").Value = "Polkadot"
$ws.Range("C15").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.96"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.07%  "
$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.875.79"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.63%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "35.581.33"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.94%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "73.81"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.11%  "
$ws.Range("E19").Value = "  -2.55%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "247.34"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.46%  "
$ws.Range("E21").Value = "  -1.17%  "
$ws.Range("E22").Value = "  -2.81%  "
$ws.Range("E23").Value = "  -0.74%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.57"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.89%  "
$ws.Range("E25").Value = "  -9.57%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "166.06"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.01%  "
$ws.Range("E27").Value = "  -1.78%  "
$ws.Range("E28").Value = "  -1.75%  "
$ws.Range("E29").Value = "  -4.04%  "
$ws.Range("E30").Value = "  -0.01%  "
$ws.Range("E31").Value = "  +7.70%  "
$ws.Range("E32").Value = "  -2.86%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0581"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.26%  "
$ws.Range("E34").Value = "  -0.20%  "
$ws.Range("E35").Value = "  -0.71%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.852"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -6.57%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.01"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.74%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.55"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -22.46%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0685"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +5.57%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "17.22"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.60%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "97.86"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.33%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0214"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.94%  "
$ws.Range("E43").Value = "  -2.36%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.299.55"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.68%  "
$ws.Range("E45").Value = "  -2.62%  "
$ws.Range("E46").Value = "  +7.67%  "
$ws.Range("E47").Value = "  -1.28%  "
$ws.Range("E48").Value = "  -0.42%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "12.15"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +4.10%  "
$ws.Range("B50").Value = "MultiversX"
$ws.Range("C50").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "43.40"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -4.14%  "
$ws.Range("B51").Value = "FraxShare"
$ws.Range("C51").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.35"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -5.70%  "
